# "Add files via upload" - append the latest SIP (2025-04-01) entries to
# the mutual fund database sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows continue straight after the last existing data row (22),
# so start by cloning that row's formatting (number formats / alignment
# styles for each column) down into the five new rows.
$ws.Range("A22:F22").Copy()
$ws.Range("A23:F27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newEntries = @(
    @{ SNo = 22; Fund = "Parag Parikh Flexi cap";          Symbol = "0P0000YWL1.BO"; Date = "2025-04-01"; Action = "SIP"; Amount = 1500 },
    @{ SNo = 23; Fund = "Nippon India Small cap";           Symbol = "0P0000XVFY.BO"; Date = "2025-04-01"; Action = "SIP"; Amount = 1000 },
    @{ SNo = 24; Fund = "DSP Nifty 50 Equal Weightage";     Symbol = "0P0001BOXZ.BO"; Date = "2025-04-01"; Action = "SIP"; Amount = 500 },
    @{ SNo = 25; Fund = "DSP Nifty Next 50";                Symbol = "0P0001FTFQ.BO"; Date = "2025-04-01"; Action = "SIP"; Amount = 500 },
    @{ SNo = 26; Fund = "Edelweiss Nifty 100 Quality 30";   Symbol = "0P0001NI59.BO"; Date = "2025-04-01"; Action = "SIP"; Amount = 500 }
)

$rowIndex = 23
foreach ($entry in $newEntries) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry.SNo
    $ws.Cells.Item($rowIndex, 2).Value = $entry.Fund
    $ws.Cells.Item($rowIndex, 3).Value = $entry.Symbol
    $ws.Cells.Item($rowIndex, 4).Value = $entry.Date
    $ws.Cells.Item($rowIndex, 5).Value = $entry.Action
    $ws.Cells.Item($rowIndex, 6).Value = $entry.Amount
    $rowIndex++
}

# Match the author's final selection/view state.
$ws.Range("F27").Select()
